$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily Orders")

# Insert a new row at row 2, shifting existing rows (2,3) down to (3,4)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new order data
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "2026-01-13 10:38"
$ws.Range("C2").Value = "Pooja"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "9096648553"
$ws.Range("E2").Value = "A 1608, Pune 411045"
$ws.Range("F2").Value = "Girl Haldi Kunku Set x1"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"
$ws.Range("J2").Value = "'"
$ws.Range("K2").Value = "'"
$ws.Range("L2").Value = "'"

# Update the Summary sheet totals
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("A2").Value = 3
$ws2.Range("B2").Value = 3
